# "Solucionado editar planificacion horaria"
# Rewrite the weekly/monthly schedule ("Planificacion") sheet:
#   - header A1 becomes "EMPLEADO" (was the last shared-string entry before)
#   - the two employee-id rows get new id numbers
#   - the schedule cells switch from the single "OFICINA" value to the new
#     "oficina"/"oficinas" (row 3) and "pasantia"/"pasantias" (row 2) values
#   - weekend columns (SABADO/DOMINGO) are left blank, same as before
#   - selection moves to A5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planificacion")

# Header row (row 1, column A) -- the per-day headers in B1:AF1 are untouched.
$ws.Range("A1").Value = "EMPLEADO"

# Employee id cells -- keep them as text (matching the original, which stored
# these as text-typed shared strings) by forcing a text entry.
$ws.Range("A2").Value = "'1724157688"
$ws.Range("A3").Value = "'1722128905"

# Work-day columns (weekends G/H, N/O, U/V, AB/AC stay empty, as in the source).
$workdayCols = @("B","C","D","E","F","I","J","K","L","M","P","Q","R","S","T","W","X","Y","Z","AA","AD","AE","AF")

foreach ($col in $workdayCols) {
    $ws.Range($col + "2").Value = "pasantia"
    $ws.Range($col + "3").Value = "oficina"
}

# A couple of cells use the plural form instead of the generic one.
$ws.Range("C2").Value = "pasantias"
$ws.Range("E3").Value = "oficinas"

$ws.Range("A5").Select()
